# Insert a new data row at sheet row 168 (pushes existing rows 168-248 down
# to 169-249) and populate it with the new weekly price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(168).Insert()

$ws.Cells.Item(168, 1).Value = 3
$ws.Cells.Item(168, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(168, 3).Value = 'Coquimbo'
$ws.Cells.Item(168, 4).Value = 44523
$ws.Cells.Item(168, 5).Value = 5
$ws.Cells.Item(168, 6).Value = 100112040
$ws.Cells.Item(168, 7).Value = 'Cilantro'
$ws.Cells.Item(168, 8).Value = 'Sin especificar'
$ws.Cells.Item(168, 9).Value = 'Primera'
$ws.Cells.Item(168, 10).Value = 160
$ws.Cells.Item(168, 11).Value = 3000
$ws.Cells.Item(168, 12).Value = 3000
$ws.Cells.Item(168, 13).Value = 3000
$ws.Cells.Item(168, 14).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(168, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(168, 16).Value = 1000
$ws.Cells.Item(168, 17).Value = 3
$ws.Cells.Item(168, 18).Value = 'Hortaliza'
